$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9407384395599365
$ws.Range("B1").Value = 1.97470235824585
$ws.Range("C1").Value = 7.547750473022461
$ws.Range("D1").Value = 2.725046157836914
$ws.Range("E1").Value = 1.230742931365967
